# M13 over ISO-25010 is vervallen.
# Maatregel M13 "Het project gebruikt ISO-25010 voor de specificatie van
# productkwaliteitseisen" is vervallen. Delete the entire row for M13
# (row 40 on the checklist sheet) and let Excel shift everything below it up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows(40).Delete()
